$d = $word.ActiveDocument

# 1) The lone paragraph in the body used the custom "Note" paragraph
#    style; drop that explicit style assignment (falls back to Normal,
#    which Word omits from the saved pPr entirely).
$p = $d.Paragraphs.Item(1)
$p.Range.ParagraphFormat.Style = "Normal"

# 2) The custom style "MarginNoteRIght" (note the stray capital I) was
#    a typo; rename it -- both its internal id and its display name --
#    to "MarginNoteRight". Word's object model only lets us edit the
#    display name (NameLocal) of an *existing* style, so recreate the
#    style under the corrected id and restore its settings.
$styleName = "MarginNoteRIght"
$fixedName = "MarginNoteRight"

$old = $d.Styles.Item($styleName)
$baseStyleName = $old.BaseStyle.NameLocal
$quick = $old.QuickStyle
$styleType = $old.Type

$old.Delete()

$new = $d.Styles.Add($fixedName, $styleType)
$new.BaseStyle = $d.Styles.Item($baseStyleName)
$new.QuickStyle = $quick
